# Apply the diff described in the commit:
#  - Metadata sheet: bump Version 5.0.0 -> 6.0.0, refresh the Date, give
#    Publisher a value, replace one of the duplicated "Contact / No display
#    for ContactDetail" rows with a "Jurisdiction / United States of
#    America" row, and delete the other duplicate row (which shifts the
#    remaining rows up by one, dropping the table from 21 to 20 rows).
#  - Elements sheet: update the Short/Definition values on row 2 to match
#    the new Title/Description text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date update
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Remove the second of the two identical "Contact" rows (row 11); this
# shifts every following row up by one, turning the 21-row table into a
# 20-row table without disturbing the (untouched) values in those rows.
$meta.Rows.Item(11).Delete()

# Turn the remaining "Contact" row (now row 10) into the new
# "Jurisdiction" row.
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# Elements sheet: row 2's Short/Definition columns (K/L) now reference the
# job-family title/description instead of the generic Extension text.
$elements = $wb.Worksheets.Item("Elements")
$elements.Cells.Item(2, 11).Value = "Employee Job Family"
$elements.Cells.Item(2, 12).Value = "Code for the job family of the employee"
